# Proofreading pass on the figure captions / alt text in appendix2.docx:
#   - "climatic scenarios" -> "climate scenarios" (Figure 1 & Figure 2)
#   - Figure 2 wording: "a change discharge as" -> "a change in discharge as"
#     and "inflating, diminishing, or remaining stable" ->
#     "inflating, remaining stable, or diminishing"
# Each change must be applied both to the visible caption paragraph text and
# to the picture's alternative text (the wp:docPr/@descr), which Word keeps
# in sync with the caption.

$d = $word.ActiveDocument

# --- Visible caption paragraphs -------------------------------------------
# Find/Replace only touches document text, so scope a separate Find to each
# caption paragraph (Find.Execute stops after the first hit, and the two
# captions each contain one "climatic").
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text

    if ($t -like "Figure 1:*") {
        $p.Range.Find.Execute("climatic", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "climate", 2) | Out-Null
    }
    elseif ($t -like "Figure 2:*") {
        $p.Range.Find.Execute("a change discharge as", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "a change in discharge as", 2) | Out-Null
        $p.Range.Find.Execute("inflating, diminishing, or remaining stable", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "inflating, remaining stable, or diminishing", 2) | Out-Null
        $p.Range.Find.Execute("climatic", $true, $false, $false, $false, $false, `
                               $true, 1, $false, "climate", 2) | Out-Null
    }
}

# --- Picture alternative text (wp:docPr/@descr) ----------------------------
# Read-modify-write via .Replace() so the existing escaped "<"/">" characters
# in the Figure 1 description are preserved byte-for-byte.
foreach ($shp in $d.InlineShapes) {
    $desc = $shp.AlternativeText

    if ($desc -like "Figure 1:*") {
        $shp.AlternativeText = $desc.Replace("climatic", "climate")
    }
    elseif ($desc -like "Figure 2:*") {
        $desc = $desc.Replace("a change discharge as", "a change in discharge as")
        $desc = $desc.Replace("inflating, diminishing, or remaining stable", "inflating, remaining stable, or diminishing")
        $desc = $desc.Replace("climatic", "climate")
        $shp.AlternativeText = $desc
    }
}
